# Revamp PricingDownload: append the newly-priced order line items beneath
# the existing rows (sheet grows from A1:E5 to A1:E8).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows, matching the existing sheet's layout: SKU | Item | Quantity | Cost Per | Total Cost
$newRows = @(
    @("33282", "Cont 12 oz - Combo (microwavable)", "2", "24.44", "48.88"),
    @("15509", "Container - Deli (32oz)",            "1", "59.94", "59.94"),
    @("14909", "Lid - Deli (6/32oz)",                 "2", "24.82", "49.64")
)

$startRow = 6
$endRow = $startRow + $newRows.Length - 1

# The source data (like the existing rows) is text, not numeric -- format the
# target range as Text first so values like "33282" aren't auto-coerced to
# numbers when assigned, matching the existing cells' string typing.
$rng = $ws.Range("A$startRow`:E$endRow")
$rng.NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
}

# Restore the default "Normal" style so the new cells carry no explicit
# formatting override, same as the pre-existing rows.
$rng.Style = "Normal"
